$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.300.99"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "3.159.68"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'606.55"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'147.77"
$ws.Range("E6").Value = "  -6.91%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.155.10"
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -4.07%  "
$ws.Range("E10").Value = "  -6.72%  "
$ws.Range("D11").Value = "'5.54"
$ws.Range("E11").Value = "  -6.79%  "
$ws.Range("D12").Value = "'0.477"
$ws.Range("E12").Value = "  -6.15%  "
$ws.Range("E13").Value = "  -7.79%  "
$ws.Range("D14").Value = "'35.74"
$ws.Range("E14").Value = "  -9.65%  "
$ws.Range("D15").Value = "3.674.37"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "64.284.73"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "3.155.79"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").Value = "'6.96"
$ws.Range("E19").Value = "  -6.05%  "
$ws.Range("D20").Value = "'481.82"
$ws.Range("E20").Value = "  -5.31%  "
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("D23").Value = "'7.74"
$ws.Range("E23").Value = "  -4.24%  "
$ws.Range("D24").Value = "'13.74"
$ws.Range("E24").Value = "  -7.66%  "
$ws.Range("D25").Value = "'83.75"
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").Value = "'8.49"
$ws.Range("E28").Value = "  -6.90%  "
$ws.Range("E29").Value = "  -8.42%  "
$ws.Range("D30").Value = "'6.78"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'0.113"
$ws.Range("E31").Value = "  -31.33%  "
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'26.27"
$ws.Range("E34").Value = "  -6.91%  "
$ws.Range("E35").Value = "  -4.66%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = "  -6.99%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'54.25"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").Value = "0.0₃0721"
$ws.Range("E38").Value = "  -10.54%  "
$ws.Range("D39").Value = "'455.27"
$ws.Range("E39").Value = "  -8.24%  "
$ws.Range("E40").Value = "  -13.83%  "
$ws.Range("D41").Value = "'0.0396"
$ws.Range("E41").Value = "  -7.68%  "
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("E43").Value = "  -8.33%  "
$ws.Range("D44").Value = "2.850.59"
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  -9.92%  "
$ws.Range("E46").Value = "  -8.46%  "
$ws.Range("D47").Value = "'26.41"
$ws.Range("E47").Value = "  -8.02%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("D51").Value = "'118.79"
$ws.Range("E51").Value = "  -2.26%  "
